$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 23 down into row 24 (keeps cell styles consistent
# with the rest of the table), then overwrite with the new row's values.
$ws.Range("A23:B23").Copy()
$ws.Range("A24:B24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D23:F23").Copy()
$ws.Range("D24:F24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Binary Tree Postorder Traversal"
$ws.Range("D24").Value = "Tree"
$ws.Range("E24").Value = "medium"
$ws.Range("F24").Value = "leetcode 145"

$ws.Range("D28").Select()
